# Weekly price-sheet update: insert a new week's record as row 43,
# pushing the existing rows 43-167 down to 44-168.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(43).Insert()

$ws.Cells.Item(43, 1).Value = 3
$ws.Cells.Item(43, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(43, 3).Value = "Coquimbo"
$ws.Cells.Item(43, 4).Value = 44497
$ws.Cells.Item(43, 5).Value = 5
$ws.Cells.Item(43, 6).Value = 100112001
$ws.Cells.Item(43, 7).Value = "Berenjena"
$ws.Cells.Item(43, 8).Value = "Sin especificar"
$ws.Cells.Item(43, 9).Value = "Primera"
$ws.Cells.Item(43, 10).Value = 90
$ws.Cells.Item(43, 11).Value = 8000
$ws.Cells.Item(43, 12).Value = 8500
$ws.Cells.Item(43, 13).Value = 8222
$ws.Cells.Item(43, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(43, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(43, 16).Value = 137
$ws.Cells.Item(43, 17).Value = 60
$ws.Cells.Item(43, 18).Value = "Hortaliza"
